$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.970.11"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +4.68%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.878.78"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.60%  "

$ws.Range("E4").Value = "  +0.05%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "279.25"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.98%  "

$ws.Range("E6").Value = "  +0.06%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5293"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +4.06%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3469"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.12%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06975"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +4.56%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.17"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.11%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.8101"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.36%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07858"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.23%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.915.16"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.66%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "90.22"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.04%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.165"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.97%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "14.56"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +4.40%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008091"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("E19").Value = "  +0.07%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "27.006.61"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +4.70%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "2.137.05"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +4.86%  "

$ws.Range("E22").Value = "  +0.78%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.05"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.78%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.197"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.66%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.351"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +7.14%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "146.48"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.51%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "17.43"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.51%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.670"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.34%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "114.38"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.59%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.360"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.78%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.339"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.71%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.08940"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.76%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.04947"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.72%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.180"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.30%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.7377"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.63%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.898"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.79%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.308"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +5.57%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.391"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +6.38%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01859"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +1.12%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5178"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.15%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.9635"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.35%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "115.93"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.64%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.206"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.14%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.119"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.88%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.07%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.4536"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.02%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.1348"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.90%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "9.396"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.41%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "36.44"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.74%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.508"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.61%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.05947"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.22%  "
